# Applies the "Added more options for creation of design matrix, and
# senstype ref." commit to design_input_example1.xlsx.

$wb  = $excel.ActiveWorkbook
$wsGeneral  = $wb.Worksheets.Item(1)   # general_input (unchanged)
$wsDesign   = $wb.Worksheets.Item(2)   # designinput
$wsDefaults = $wb.Worksheets.Item(3)   # defaultvalues

# ---------------------------------------------------------------------
# 1. designinput sheet: header / label text updates
# ---------------------------------------------------------------------
$wsDesign.Range("E1").Value = "casename1"
$wsDesign.Range("G1").Value = "casename2"

$wsDesign.Range("D3").Value = "FAULT_POSITION"
$wsDesign.Range("D4").Value = "DC_MODEL"
$wsDesign.Range("D5").Value = "OWC1"
$wsDesign.Range("D6").Value = "OWC2"
$wsDesign.Range("D7").Value = "OWC3"
$wsDesign.Range("D8").Value = "MULTZ_ILE"

# D2 ("RMS_SEED") is cleared out entirely (C2 keeps its "seed" text).
$wsDesign.Range("D2").Value = ""

# ---------------------------------------------------------------------
# 2. defaultvalues sheet: label text updates (uppercase the identifiers)
# ---------------------------------------------------------------------
$wsDefaults.Range("A3").Value  = "FAULT_POSITION"
$wsDefaults.Range("A4").Value  = "DC_MODEL"
$wsDefaults.Range("A5").Value  = "OWC1"
$wsDefaults.Range("A6").Value  = "OWC2"
$wsDefaults.Range("A7").Value  = "OWC3"
$wsDefaults.Range("A8").Value  = "MULTZ_ILE"
$wsDefaults.Range("A9").Value  = "PARAM1"
$wsDefaults.Range("A10").Value = "PARAM2"
$wsDefaults.Range("A11").Value = "PARAM3"
$wsDefaults.Range("A12").Value = "PARAM4"

# ---------------------------------------------------------------------
# 3. Comment text updates on designinput (E1 / I1); B1 stays the same.
# ---------------------------------------------------------------------
$commentE1 = $wsDesign.Range("E1").Comment
$commentE1.Text("For scenario sensitivities provide names for case1 and case2 and values. Values can be string or numbers")

$openQuote  = [char]0x201C
$closeQuote = [char]0x201D
$enDash     = [char]0x2013

$commentI1Text = "Distname and dist_param1, .. only for sensitivities of type " + $openQuote + "dist" + $closeQuote + ". " + "`n" + `
  "The order of distribution parameters is predefined: " + "`n" + `
  "normal(mean, std dev,min, max)     " + $enDash + " where min/max is optional and will give truncated gaussian" + "`n" + `
  "lognormal(mean, stddev) " + "`n" + `
  "uniform(from,to)" + "`n" + `
  "loguniform(from, to)" + "`n" + `
  "triangular(low, mode, high)" + "`n" + `
  "discrete((value1, value2, value3,..,value_n) (weight1, weight2, weight3,..weight_n)). Discrete uniform if no weights are given"

$commentI1 = $wsDesign.Range("I1").Comment
$commentI1.Text($commentI1Text)

# ---------------------------------------------------------------------
# 4. Style updates: E3:E5 alignment changes to "left".
# ---------------------------------------------------------------------
$wsDesign.Range("E3:E5").HorizontalAlignment = -4131   # xlHAlignLeft

# Column D widens.
$wsDesign.Columns.Item(4).ColumnWidth = 15

# ---------------------------------------------------------------------
# 5. View / selection state: designinput becomes the active sheet with
#    selection Q35; defaultvalues is no longer the active tab and its
#    selection moves to B12.
# ---------------------------------------------------------------------
$wsDefaults.Range("B12").Select()
$wsDesign.Activate()
$wsDesign.Range("Q35").Select()
